$wb = $excel.ActiveWorkbook

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 370.2857
$ws.Range("I39").Value = 265.33334
$ws.Range("K39").Value = 796.0000200000001
$ws.Range("M39").Value = -500.0000200000001

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4292.9
$ws.Range("I100").Value = 4658.778
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 4658.778
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -4117.778
$ws.Range("N100").Value = -2082

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 326.25
$ws.Range("I111").Value = 326.25
$ws.Range("K111").Value = 978.75
$ws.Range("M111").Value = 2088.25

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 624.75

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1304.25

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 419.5
$ws.Range("I135").Value = 460
$ws.Range("J135").Value = 176.5
$ws.Range("K135").Value = 4140
$ws.Range("L135").Value = 1588.5
$ws.Range("M135").Value = -1605
$ws.Range("N135").Value = -6658.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3666.875
$ws.Range("I138").Value = 2100
$ws.Range("J138").Value = 3957.037
$ws.Range("K138").Value = 6300
$ws.Range("L138").Value = 11871.111
$ws.Range("M138").Value = -1160
$ws.Range("N138").Value = -22151.111

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18415.334
$ws.Range("I32").Value = 18415.334
$ws.Range("K32").Value = 18415.334
$ws.Range("M32").Value = -18128.334

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2036.6923
$ws.Range("I61").Value = 2023.0834
$ws.Range("K61").Value = 2023.0834
$ws.Range("M61").Value = -1811.0834

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13284.091
$ws.Range("I74").Value = 13284.091
$ws.Range("K74").Value = 13284.091
$ws.Range("M74").Value = -12410.091

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 13284.091
$ws.Range("I77").Value = 13284.091
$ws.Range("K77").Value = 66420.455
$ws.Range("M77").Value = -62052.455

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2403.6365
$ws.Range("I102").Value = 2619
$ws.Range("J102").Value = 250
$ws.Range("K102").Value = 2619
$ws.Range("L102").Value = 250
$ws.Range("M102").Value = -997
$ws.Range("N102").Value = -3494

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2361
$ws.Range("I122").Value = 2361
$ws.Range("K122").Value = 7083
$ws.Range("M122").Value = -4633

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2751.25
$ws.Range("I132").Value = 4069.8
$ws.Range("K132").Value = 12209.4
$ws.Range("M132").Value = -9679.400000000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2036.6923
$ws.Range("I136").Value = 2023.0834
$ws.Range("K136").Value = 6069.2502
$ws.Range("M136").Value = -3519.2502

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1162.5
$ws.Range("I99").Value = 1196.4
$ws.Range("J99").Value = 993
$ws.Range("K99").Value = 1196.4
$ws.Range("L99").Value = 993
$ws.Range("M99").Value = 301.5999999999999
$ws.Range("N99").Value = -3989

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1599.8
$ws.Range("I107").Value = 1599.8
$ws.Range("K107").Value = 1599.8
$ws.Range("M107").Value = 320.2

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2451.0588
$ws.Range("I134").Value = 2389
$ws.Range("K134").Value = 7167
$ws.Range("M134").Value = -4632

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2689
$ws.Range("I31").Value = 2797.8
$ws.Range("J31").Value = 2553
$ws.Range("K31").Value = 2797.8
$ws.Range("L31").Value = 2553
$ws.Range("M31").Value = -2502.8
$ws.Range("N31").Value = -3143

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2689
$ws.Range("I34").Value = 2797.8
$ws.Range("J34").Value = 2553
$ws.Range("K34").Value = 2797.8
$ws.Range("L34").Value = 2553
$ws.Range("M34").Value = -2595.8
$ws.Range("N34").Value = -2957

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6648.6665
$ws.Range("I58").Value = 2264
$ws.Range("J58").Value = 11033.333
$ws.Range("K58").Value = 2264
$ws.Range("L58").Value = 11033.333
$ws.Range("M58").Value = -2061
$ws.Range("N58").Value = -11439.333

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2207.7856
$ws.Range("I107").Value = 2264.9092
$ws.Range("K107").Value = 2264.9092
$ws.Range("M107").Value = -344.9092000000001

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4899.625
$ws.Range("I122").Value = 5466.3335
$ws.Range("J122").Value = 3199.5
$ws.Range("K122").Value = 16399.0005
$ws.Range("L122").Value = 9598.5
$ws.Range("M122").Value = -13949.0005
$ws.Range("N122").Value = -14498.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3240.2727
$ws.Range("I132").Value = 2275
$ws.Range("J132").Value = 4398.6
$ws.Range("K132").Value = 6825
$ws.Range("L132").Value = 13195.8
$ws.Range("M132").Value = -4295
$ws.Range("N132").Value = -18255.8

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6648.6665
$ws.Range("I136").Value = 2264
$ws.Range("J136").Value = 11033.333
$ws.Range("K136").Value = 6792
$ws.Range("L136").Value = 33099.999
$ws.Range("M136").Value = -4242
$ws.Range("N136").Value = -38199.999

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 267.85715
$ws.Range("I12").Value = 550.5
$ws.Range("J12").Value = 154.8
$ws.Range("K12").Value = 1651.5
$ws.Range("L12").Value = 464.4
$ws.Range("M12").Value = -1478.5
$ws.Range("N12").Value = -810.4000000000001

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 778.5
$ws.Range("I14").Value = 778.5
$ws.Range("K14").Value = 2335.5
$ws.Range("M14").Value = -2162.5

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4195
$ws.Range("J34").Value = 4195
$ws.Range("L34").Value = 12585
$ws.Range("N34").Value = -12753

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 8228.75
$ws.Range("J55").Value = 14495
$ws.Range("L55").Value = 43485
$ws.Range("N55").Value = -43839

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1305.0769
$ws.Range("I129").Value = 1435.4
$ws.Range("K129").Value = 4306.200000000001
$ws.Range("M129").Value = 693.7999999999993

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2561.5
$ws.Range("I137").Value = 2945.7144
$ws.Range("K137").Value = 8837.143199999999
$ws.Range("M137").Value = -3737.143199999999

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1172.7273
$ws.Range("I107").Value = 816.6667
$ws.Range("K107").Value = 2450.0001
$ws.Range("M107").Value = -530.0001000000002

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1772.6897
$ws.Range("I132").Value = 1539.6957
$ws.Range("K132").Value = 4619.0871
$ws.Range("M132").Value = -2089.0871

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7509.273
$ws.Range("I136").Value = 6760.2
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 20280.6
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -17730.6
$ws.Range("N136").Value = -50100

Write-Output "Applied all Sophia_Profits updates"